# "Added Sum for block"
# Fills in the previously-blank "Reduction Sum K3" (column block N:Q) and
# "Added Tiling to K3 ..." (column block R:U) result columns on Sheet2 with
# the same kind of data already present for the other version blocks
# (B:E, F:I, J:M), and extends the SUM() total row (16) to cover the new
# N:Q block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Seed / Converge after / Silhouette score / Total time (printed) ----
# (rows 5-8, merged column blocks N:Q and R:U)
$ws.Cells.Item(5, 14).Value = 42          # N5  Seed
$ws.Cells.Item(5, 18).Value = 42          # R5  Seed

$ws.Cells.Item(6, 14).Value = 32          # N6  Converge after
$ws.Cells.Item(6, 18).Value = 13          # R6  Converge after

$ws.Cells.Item(7, 14).Value = 0.61830700000000005   # N7  Silhouette Score
$ws.Cells.Item(7, 18).Value = 0.61830700000000005   # R7  Silhouette Score

$ws.Cells.Item(8, 14).Value = "1.827291 sec"         # N8  Total Time (Printed)
$ws.Cells.Item(8, 18).Value = "1.931982 sec"         # R8  Total Time (Printed)

# ---- Profiler breakdown rows 10-15 for the two new blocks ----
# (plain decimal literals - PowerShell's parser here rejects the
# "1.23E-4" scientific-notation form, but a plain decimal round-trips to
# the exact same IEEE-754 double / OOXML <v> serialization)
# Block N:Q ("Reduction Sum K3")
$ws.Cells.Item(10, 14).Value = 0.9976
$ws.Cells.Item(10, 15).Value = "1.65274s"
$ws.Cells.Item(10, 16).Value = 1
$ws.Cells.Item(10, 17).Value = "1.65274s "

$ws.Cells.Item(11, 14).Value = 0.0017
$ws.Cells.Item(11, 15).Value = "2.7870ms"
$ws.Cells.Item(11, 16).Value = 13
$ws.Cells.Item(11, 17).Value = " 214.39us"

$ws.Cells.Item(12, 14).Value = 0.0004
$ws.Cells.Item(12, 15).Value = "721.10us"
$ws.Cells.Item(12, 16).Value = 44
$ws.Cells.Item(12, 17).Value = "16.388us"

$ws.Cells.Item(13, 14).Value = 0.0002
$ws.Cells.Item(13, 15).Value = "362.84us"
$ws.Cells.Item(13, 16).Value = 45
$ws.Cells.Item(13, 17).Value = "8.0630us"

$ws.Cells.Item(14, 14).Value = 0.0001
$ws.Cells.Item(14, 15).Value = "133.02us"
$ws.Cells.Item(14, 16).Value = 28
$ws.Cells.Item(14, 17).Value = "4.7500us"

$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = "15.935us"
$ws.Cells.Item(15, 16).Value = 26
$ws.Cells.Item(15, 17).Value = "612ns "

# Block R:U ("Added Tiling to K3 ...")
$ws.Cells.Item(10, 18).Value = 0.9961
$ws.Cells.Item(10, 19).Value = "1.92609s"
$ws.Cells.Item(10, 20).Value = 1
$ws.Cells.Item(10, 21).Value = "1.92609s "

$ws.Cells.Item(11, 18).Value = 0.0026
$ws.Cells.Item(11, 19).Value = " 4.9877ms"
$ws.Cells.Item(11, 20).Value = 13
$ws.Cells.Item(11, 21).Value = "383.67us"

$ws.Cells.Item(12, 18).Value = 0.001
$ws.Cells.Item(12, 19).Value = "1.8528ms"
$ws.Cells.Item(12, 20).Value = 44
$ws.Cells.Item(12, 21).Value = "42.108us"

$ws.Cells.Item(13, 18).Value = 0.0002
$ws.Cells.Item(13, 19).Value = "380.00us"
$ws.Cells.Item(13, 20).Value = 45
$ws.Cells.Item(13, 21).Value = "8.4440us"

$ws.Cells.Item(14, 18).Value = 0.0001
$ws.Cells.Item(14, 19).Value = "238.15us"
$ws.Cells.Item(14, 20).Value = 28
$ws.Cells.Item(14, 21).Value = "8.5050us"

$ws.Cells.Item(15, 18).Value = 0
$ws.Cells.Item(15, 19).Value = "26.336us"
$ws.Cells.Item(15, 20).Value = 26
$ws.Cells.Item(15, 21).Value = "1.0120us"

# ---- Totals row 16: extend the SUM() block to N:Q, matching F:I / J:M ----
$ws.Range("F16").Copy()
$ws.Range("N16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G16").Copy()
$ws.Range("O16").PasteSpecial(-4122)
$ws.Range("H16").Copy()
$ws.Range("P16").PasteSpecial(-4122)
$ws.Range("I16").Copy()
$ws.Range("Q16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("N16").Formula = "=SUM(N10:N15)"
$ws.Range("O16").Value = "1.65675s"
$ws.Range("Q16").Value = "1.652984s"
